# 1st testcase of forgotPassword scenario added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Test Suite" sheet

# A4 is written before B3 so the shared-string table grows in the same
# order the original author's Excel session produced it in.
$ws.Range("A4").Value = "Forgot Password"

# Row 3 ("Login") had an incorrect description left over from a
# copy/paste of the Catalogue row - fix it to describe the Login suite.
$ws.Range("B3").Value = "Login suite description"

# Finish populating the new "Forgot Password" test-suite row.
$ws.Range("B4").Value = "Forgot Password suite description"
$ws.Range("C4").Value = "YES"

# Leave the selection where the author's session ended up.
$ws.Range("C7").Select()
